$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25: quantity bump + new line-total formula ---
$ws.Range("B25").Value = 5
$ws.Range("E25").Formula = "=B25*D25"

# --- Row 26: new line-total formula ---
$ws.Range("E26").Formula = "=B26*D26"

# --- Rows 50-54: components from other suppliers / JCU storeroom ---
$ws.Range("A50").Value = "1N4001 "
$ws.Range("B50").Value = 2
$ws.Range("C50").Value = "Diode 1A 50V"

$ws.Range("B51").Value = 2
$ws.Range("C51").Value = "4k7 Resistor"

$ws.Range("B52").Value = 1
$ws.Range("C52").Value = "0.33uF Capacitor"

$ws.Range("B53").Value = 1
$ws.Range("C53").Value = "0.1 uF Capacitor"

$ws.Range("B54").Value = 1
$ws.Range("C54").Value = "1uF Capacitor"

# --- View state: scrolled/selected further down the sheet, zoomed to 115% ---
$ws.Range("A52").Select()
$excel.ActiveWindow.Zoom = 115
